# Update LR-pair matrix (NATMI output) per Dr Hou's advice:
# ligand/receptor-expressing-cell counts (E,K columns) and all downstream
# weighted/specificity statistics were recomputed, and the sending x target
# cluster matrix now also includes M2 as a target cluster (rows 4,8,12,16,20),
# expanding the sheet from 15 to 20 data rows (A2:T16 -> A2:T21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is one full data row (columns A..T) for sheet rows 2..21:
# Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
# Ligand-expressing cells, Ligand detection rate, Ligand average expression
# value, Ligand total expression value, Ligand derived specificity of
# average/total expression value, Receptor-expressing cells, Receptor
# detection rate, Receptor average/total expression value, Receptor derived
# specificity of average/total expression value, Edge average/total
# expression weight, Edge average/total expression derived specificity.
$data = @(
    @('ECs', 'Pdgfa', 'Pdgfra', 'ECs', 3, 1, 9.327880333333333, 27.983641, 0.2334346564853683, 0.2334346564853683, 3, 1, 2.891565, 8.674695, 0.02565374714703297, 0.02565374714703297, 26.972172296055, 242.749550664495, 0.005988473652830137, 0.005988473652830138),
    @('ECs', 'Pdgfa', 'Pdgfra', 'FAPs', 3, 1, 9.327880333333333, 27.983641, 0.2334346564853683, 0.2334346564853683, 3, 1, 109.4279763333333, 328.283929, 0.9708367737425379, 0.970836773742538, 1020.731068356165, 9186.579615205488, 0.2266269487819525, 0.2266269487819526),
    @('ECs', 'Pdgfa', 'Pdgfra', 'M2', 3, 1, 9.327880333333333, 27.983641, 0.2334346564853683, 0.2334346564853683, 1, 0.3333333333333333, 0.06646400000000001, 0.199392, 0.00058966360790105, 0.00058966360790105, 0.6199682384746668, 5.579714146272, 0.0001376479217523045, 0.0001376479217523045),
    @('ECs', 'Pdgfa', 'Pdgfra', 'sCs', 3, 1, 9.327880333333333, 27.983641, 0.2334346564853683, 0.2334346564853683, 3, 1, 0.3291073333333334, 0.987322, 0.002919815502528088, 0.002919815502528088, 3.069873822155778, 27.628864399402, 0.0006815861288332971, 0.0006815861288332971),
    @('FAPs', 'Pdgfa', 'Pdgfra', 'ECs', 3, 1, 7.038387, 21.115161, 0.1761389933021313, 0.1761389933021313, 3, 1, 2.891565, 8.674695, 0.02565374714703297, 0.02565374714703297, 20.351953505655, 183.167581550895, 0.004518625196905809, 0.004518625196905809),
    @('FAPs', 'Pdgfa', 'Pdgfra', 'FAPs', 3, 1, 7.038387, 21.115161, 0.1761389933021313, 0.1761389933021313, 3, 1, 109.4279763333333, 328.283929, 0.9708367737425379, 0.970836773742538, 770.196446060841, 6931.768014547569, 0.1710022119876996, 0.1710022119876996),
    @('FAPs', 'Pdgfa', 'Pdgfra', 'M2', 3, 1, 7.038387, 21.115161, 0.1761389933021313, 0.1761389933021313, 1, 0.3333333333333333, 0.06646400000000001, 0.199392, 0.00058966360790105, 0.00058966360790105, 0.4677993535680001, 4.210194182112001, 0.0001038627542825936, 0.0001038627542825936),
    @('FAPs', 'Pdgfa', 'Pdgfra', 'sCs', 3, 1, 7.038387, 21.115161, 0.1761389933021313, 0.1761389933021313, 3, 1, 0.3291073333333334, 0.987322, 0.002919815502528088, 0.002919815502528088, 2.316384776538, 20.847462988842, 0.0005142933632432539, 0.0005142933632432539),
    @('M1', 'Pdgfa', 'Pdgfra', 'ECs', 3, 1, 5.440503, 16.321509, 0.1361511837125786, 0.1361511837125786, 3, 1, 2.891565, 8.674695, 0.02565374714703297, 0.02565374714703297, 15.731568057195, 141.584112514755, 0.003492788040731725, 0.003492788040731725),
    @('M1', 'Pdgfa', 'Pdgfra', 'FAPs', 3, 1, 5.440503, 16.321509, 0.1361511837125786, 0.1361511837125786, 3, 1, 109.4279763333333, 328.283929, 0.9708367737425379, 0.9708367737425379, 595.343233525429, 5358.08910172886, 0.1321805759367474, 0.1321805759367474),
    @('M1', 'Pdgfa', 'Pdgfra', 'M2', 3, 1, 5.440503, 16.321509, 0.1361511837125786, 0.1361511837125786, 1, 0.3333333333333333, 0.06646400000000001, 0.199392, 0.00058966360790105, 0.00058966360790105, 0.361597591392, 3.254378322528, [double]"8.028339820795775E-05", [double]"8.028339820795775E-05"),
    @('M1', 'Pdgfa', 'Pdgfra', 'sCs', 3, 1, 5.440503, 16.321509, 0.1361511837125786, 0.1361511837125786, 3, 1, 0.3291073333333334, 0.987322, 0.002919815502528088, 0.002919815502528088, 1.790509434322, 16.114584908898, 0.0003975363368915367, 0.0003975363368915367),
    @('M2', 'Pdgfa', 'Pdgfra', 'ECs', 3, 1, 4.327482, 12.982446, 0.1082973020683707, 0.1082973020683707, 3, 1, 2.891565, 8.674695, 0.02565374714703297, 0.02565374714703297, 12.51319548933, 112.61875940397, 0.002778231603967833, 0.002778231603967833),
    @('M2', 'Pdgfa', 'Pdgfra', 'FAPs', 3, 1, 4.327482, 12.982446, 0.1082973020683707, 0.1082973020683707, 3, 1, 109.4279763333333, 328.283929, 0.9708367737425379, 0.970836773742538, 473.547597878926, 4261.928380910334, 0.1051390033450781, 0.1051390033450781),
    @('M2', 'Pdgfa', 'Pdgfra', 'M2', 3, 1, 4.327482, 12.982446, 0.1082973020683707, 0.1082973020683707, 1, 0.3333333333333333, 0.06646400000000001, 0.199392, 0.00058966360790105, 0.00058966360790105, 0.287621763648, 2.588595872832, [double]"6.385897786358531E-05", [double]"6.385897786358531E-05"),
    @('M2', 'Pdgfa', 'Pdgfra', 'sCs', 3, 1, 4.327482, 12.982446, 0.1082973020683707, 0.1082973020683707, 3, 1, 0.3291073333333334, 0.987322, 0.002919815502528088, 0.002919815502528088, 1.424206061068, 12.817854549612, 0.000316208141461196, 0.000316208141461196),
    @('sCs', 'Pdgfa', 'Pdgfra', 'ECs', 3, 1, 13.82502566666667, 41.475077, 0.3459778644315512, 0.3459778644315512, 3, 1, 2.891565, 8.674695, 0.02565374714703297, 0.02565374714703297, 39.975960341835, 359.783643076515, 0.008875628652597467, 0.008875628652597467),
    @('sCs', 'Pdgfa', 'Pdgfra', 'FAPs', 3, 1, 13.82502566666667, 41.475077, 0.3459778644315512, 0.3459778644315512, 3, 1, 109.4279763333333, 328.283929, 0.9708367737425379, 0.970836773742538, 1512.844581459726, 13615.60123313753, 0.3358880336910603, 0.3358880336910604),
    @('sCs', 'Pdgfa', 'Pdgfra', 'M2', 3, 1, 13.82502566666667, 41.475077, 0.3459778644315512, 0.3459778644315512, 1, 0.3333333333333333, 0.06646400000000001, 0.199392, 0.00058966360790105, 0.00058966360790105, 0.9188665059093335, 8.269798553184, 0.0002040105557946088, 0.0002040105557946088),
    @('sCs', 'Pdgfa', 'Pdgfra', 'sCs', 3, 1, 13.82502566666667, 41.475077, 0.3459778644315512, 0.3459778644315512, 3, 1, 0.3291073333333334, 0.987322, 0.002919815502528088, 0.002919815502528088, 4.549917330421556, 40.949255973794, 0.001010191532098804, 0.001010191532098804)
)

$r = 2
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}
